$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("C21").Value = [double]"5.77690575018865e-294"
$ws.Range("C22").Value = [double]"8.73829380307184e-224"
$ws.Range("C23").Value = [double]"2.446858791308649e-181"
$ws.Range("C24").Value = [double]"7.905047162642377e-153"
$ws.Range("C25").Value = [double]"2.541150349171852e-132"
$ws.Range("C26").Value = [double]"7.739005551381953e-117"
$ws.Range("C27").Value = [double]"1.02168934875685e-104"
$ws.Range("C28").Value = [double]"5.839870991935935e-95"
$ws.Range("C29").Value = [double]"6.279507388237319e-87"
$ws.Range("C30").Value = [double]"3.392151373416583e-80"
$ws.Range("C31").Value = [double]"1.820674350880554e-74"
$ws.Range("C32").Value = [double]"1.582495736526903e-69"
$ws.Range("C33").Value = [double]"3.190064581093582e-65"
$ws.Range("C34").Value = [double]"1.953913053992746e-61"
$ws.Range("C35").Value = [double]"4.497587531792497e-55"
$ws.Range("C36").Value = [double]"6.264054735439408e-50"
$ws.Range("C37").Value = [double]"1.655699656302115e-40"
$ws.Range("C38").Value = [double]"5.661966877745586e-34"
$ws.Range("C39").Value = [double]"1.120799781879818e-28"
$ws.Range("C40").Value = [double]"2.419149930997881e-24"
$ws.Range("C41").Value = [double]"6.303221746413665e-21"
$ws.Range("C42").Value = [double]"3.411646885679765e-18"
$ws.Range("C43").Value = [double]"4.191793247700701e-14"
$ws.Range("C44").Value = [double]"3.401799724614606e-11"
$ws.Range("C45").Value = [double]"5.080053103998901e-09"
$ws.Range("C46").Value = [double]"2.465654323714587e-07"
$ws.Range("C47").Value = [double]"5.468177354599642e-06"
$ws.Range("C48").Value = [double]"0.00143404764608668"
$ws.Range("C49").Value = [double]"0.05920742882708949"
$ws.Range("C50").Value = [double]"0.8627231011645121"
$ws.Range("C51").Value = [double]"6.612629387137179"
$ws.Range("C52").Value = [double]"123.51503246558"
$ws.Range("C53").Value = [double]"940.8579914002082"
$ws.Range("C54").Value = [double]"4221.813601361049"
$ws.Range("C55").Value = [double]"13401.75826439215"
$ws.Range("C56").Value = [double]"69911.28735052758"
$ws.Range("C57").Value = [double]"212400.8247670185"
$ws.Range("C58").Value = [double]"466958.7204594503"
$ws.Range("C59").Value = [double]"834136.9215335664"
$ws.Range("C60").Value = [double]"1294931.889329147"
$ws.Range("C61").Value = [double]"1821184.720867028"
